# Update the dSF (column F) values for the eovaldi_nathan sheet.
# This mirrors a "repull data, push all data, mean calculation" refresh
# where only the dSF column values changed for most data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = -1
    4  = 2
    5  = -4
    6  = 6
    7  = 2
    8  = 1
    9  = 1
    10 = -3
    11 = -1
    12 = 4
    13 = 2
    14 = -5
    15 = 0
    16 = 9
    17 = 7
    19 = -1
    20 = -4
    21 = 5
    23 = -2
    24 = 3
    26 = -1
    27 = -3
    28 = 1
    29 = 1
    30 = -1
    31 = 1
    32 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
